# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns,
# and fix the Algorand/Stellar row ordering (rows 37-38).
#
# Note: several "Price" values look numeric (e.g. "1.001", "0.9995") but
# must stay as literal text to match the original cell formatting/type.
# Prefixing with a leading apostrophe forces Excel to treat the value as
# text instead of auto-converting it to a number, and resetting the
# cell Style back to "Normal" afterwards avoids leaving a stray
# quote-prefix/text style applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.351.01"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.627.01"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'0.9995"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'302.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("D7").Value = "'0.3751"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "'0.3625"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'51.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("D10").Value = "'0.08158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").Value = "'1.216"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "'22.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("D14").Value = "'6.456"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").Value = "'7.277"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "1.619.65"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "'94.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").Value = "'0.06939"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "'17.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.23%  "

$ws.Range("D21").Value = "'6.543"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'12.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.17%  "

$ws.Range("D24").Value = "23.352.74"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").Value = "'2.492"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.25%  "

$ws.Range("D26").Value = "'3.057"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").Value = "'21.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").Value = "'149.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").Value = "'5.270"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "'132.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").Value = "1.762.60"
$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("D32").Value = "'6.588"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "

$ws.Range("E33").Value = "  -5.62%  "

$ws.Range("D34").Value = "'1.057"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.16%  "

$ws.Range("E35").Value = "  +7.73%  "

$ws.Range("D36").Value = "'0.02750"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2485"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.08750"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").Value = "'5.958"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("D41").Value = "'0.6964"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.14%  "

$ws.Range("D42").Value = "'1.326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.26%  "

$ws.Range("D43").Value = "'15.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").Value = "'11.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("D45").Value = "'0.6428"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("D46").Value = "'0.9993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'2.266"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.69%  "

$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("D49").Value = "'0.07968"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "'126.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").Value = "'1.188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
